$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8274975676880558
$ws.Range("C2").Value = 0.2026139200537216
$ws.Range("E2").Value = 0.1087810853287143
$ws.Range("F2").Value = 0.4443680307746263
$ws.Range("G2").Value = 0.002426634525031254
$ws.Range("I2").Value = 0.6385094511219727
$ws.Range("L2").Value = 0.203175727809878
$ws.Range("M2").Value = 0.1985441872694409
$ws.Range("N2").Value = 1.293369678917308
$ws.Range("O2").Value = 2.417229885097413
$ws.Range("B3").Value = 0.7512529339883827
$ws.Range("C3").Value = 0.1916004488447243
$ws.Range("E3").Value = 0.1096980128860374
$ws.Range("F3").Value = 0.387822817061874
$ws.Range("G3").Value = 0.002429146138284032
$ws.Range("I3").Value = 0.6455826104999822
$ws.Range("L3").Value = 0.2004646143127573
$ws.Range("M3").Value = 0.1860604183255745
$ws.Range("N3").Value = 1.302063999176454
$ws.Range("O3").Value = 2.427600879706233
$ws.Range("B4").Value = 0.7045184088714791
$ws.Range("C4").Value = 0.1847795741160922
$ws.Range("E4").Value = 0.1102965299168699
$ws.Range("F4").Value = 0.3531389305168915
$ws.Range("G4").Value = 0.002430771511706303
$ws.Range("I4").Value = 0.6503059831908793
$ws.Range("L4").Value = 0.1988999989937028
$ws.Range("M4").Value = 0.1784516546754773
$ws.Range("N4").Value = 1.307890540282827
$ws.Range("O4").Value = 2.435587458445767
$ws.Range("B5").Value = 0.6854951171162895
$ws.Range("C5").Value = 0.1819854238413825
$ws.Range("E5").Value = 0.1105493743245409
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002431454855690387
$ws.Range("I5").Value = 0.6523263944119009
$ws.Range("L5").Value = 0.1982876215282161
$ws.Range("M5").Value = 0.1753653973838603
$ws.Range("N5").Value = 1.310387829895291
$ws.Range("O5").Value = 2.439248778321513
$ws.Range("B6").Value = 0.6823376450221872
$ws.Range("C6").Value = 0.1815205802414823
$ws.Range("E6").Value = 0.1105918994889628
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002431569593995919
$ws.Range("I6").Value = 0.6526676542480985
$ws.Range("L6").Value = 0.1981874616234407
$ws.Range("M6").Value = 0.1748538012519134
$ws.Range("N6").Value = 1.310809932378774
$ws.Range("O6").Value = 2.439881293404127
$ws.Range("B7").Value = 0.704261765717149
$ws.Range("C7").Value = 0.1847419501135903
$ws.Range("E7").Value = 0.1102999036354189
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.002430780642464943
$ws.Range("I7").Value = 0.6503328441715368
$ws.Range("L7").Value = 0.1988916380717001
$ws.Range("M7").Value = 0.1784099738407789
$ws.Range("N7").Value = 1.307923721654852
$ws.Range("O7").Value = 2.435635189905383
$ws.Range("B8").Value = 0.8011925676100304
$ws.Range("C8").Value = 0.198828707928044
$ws.Range("E8").Value = 0.1090898771582713
$ws.Range("F8").Value = 0.4248636149813336
$ws.Range("G8").Value = 0.002427483291439807
$ws.Range("I8").Value = 0.6408692822157356
$ws.Range("L8").Value = 0.2022202146972134
$ws.Range("M8").Value = 0.1942282180910908
$ws.Range("N8").Value = 1.296266254751025
$ws.Range("O8").Value = 2.420469634731063
$ws.Range("B9").Value = 0.9918608934855797
$ws.Range("C9").Value = 0.2259833581679231
$ws.Range("E9").Value = 0.1069982920322254
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.002421674798611098
$ws.Range("I9").Value = 0.625332577035838
$ws.Range("L9").Value = 0.2095388465632979
$ws.Range("M9").Value = 0.2256869047645154
$ws.Range("N9").Value = 1.277272558767564
$ws.Range("O9").Value = 2.403590992233148
$ws.Range("B10").Value = 1.132253242147726
$ws.Range("C10").Value = 0.2456430196029089
$ws.Range("E10").Value = 0.1056322804678764
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002417804250824028
$ws.Range("I10").Value = 0.6157635399931927
$ws.Range("L10").Value = 0.2153959790571349
$ws.Range("M10").Value = 0.2490593723451298
$ws.Range("N10").Value = 1.26566601191869
$ws.Range("O10").Value = 2.399057532363742
$ws.Range("B11").Value = 1.196178622355717
$ws.Range("C11").Value = 0.2545226664589393
$ws.Range("E11").Value = 0.1050477374489445
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002416128794583092
$ws.Range("I11").Value = 0.6118119837818128
$ws.Range("L11").Value = 0.2181643680961685
$ws.Range("M11").Value = 0.2597469823635308
$ws.Range("N11").Value = 1.260893956070539
$ws.Range("O11").Value = 2.398709456904811
$ws.Range("B12").Value = 1.220393029426759
$ws.Range("C12").Value = 0.2578758916685615
$ws.Range("E12").Value = 0.1048316746376452
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.002415506542154874
$ws.Range("I12").Value = 0.6103734347964149
$ws.Range("L12").Value = 0.219227579119206
$ws.Range("M12").Value = 0.2638018788175955
$ws.Range("N12").Value = 1.259159787759003
$ws.Range("O12").Value = 2.398824580735038
$ws.Range("B13").Value = 1.215177723353634
$ws.Range("C13").Value = 0.2571541308872156
$ws.Range("E13").Value = 0.1048779724698878
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.00241564001333815
$ws.Range("I13").Value = 0.6106806790537114
$ws.Range("L13").Value = 0.2189979367608146
$ws.Range("M13").Value = 0.2629282441560647
$ws.Range("N13").Value = 1.259530031274224
$ws.Range("O13").Value = 2.398788797549173
$ws.Range("B14").Value = 1.198170618868687
$ws.Range("C14").Value = 0.2547987256699003
$ws.Range("E14").Value = 0.1050298558491836
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.002416077357102172
$ws.Range("I14").Value = 0.6116924741736689
$ws.Range("L14").Value = 0.2182515412223438
$ws.Range("M14").Value = 0.2600804275203075
$ws.Range("N14").Value = 1.260749824524716
$ws.Range("O14").Value = 2.398713976862894
$ws.Range("B15").Value = 1.187754179240642
$ws.Range("C15").Value = 0.2533547556920155
$ws.Range("E15").Value = 0.1051235775505166
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002416346830921794
$ws.Range("I15").Value = 0.6123197604460984
$ws.Range("L15").Value = 0.217796288549934
$ws.Range("M15").Value = 0.2583370575230219
$ws.Range("N15").Value = 1.261506474140148
$ws.Range("O15").Value = 2.398700317647979
$ws.Range("B16").Value = 1.128076659515955
$ws.Range("C16").Value = 0.2450614201785584
$ws.Range("E16").Value = 0.1056712226543111
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002417915457038833
$ws.Range("I16").Value = 0.6160298694909905
$ws.Range("L16").Value = 0.2152171444485731
$ws.Range("M16").Value = 0.2483620061749434
$ws.Range("N16").Value = 1.26598808543266
$ws.Range("O16").Value = 2.399114811172808
$ws.Range("B17").Value = 1.091480852753591
$ws.Range("C17").Value = 0.2399573200645477
$ws.Range("E17").Value = 0.1060166195389227
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.002418899559327077
$ws.Range("I17").Value = 0.6184087868985166
$ws.Range("L17").Value = 0.2136615004606597
$ws.Range("M17").Value = 0.2422566500025454
$ws.Range("N17").Value = 1.268867388351211
$ws.Range("O17").Value = 2.399808448581751
$ws.Range("B18").Value = 1.070437662347956
$ws.Range("C18").Value = 0.2370155965576259
$ws.Range("E18").Value = 0.1062187529492311
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002419473618898123
$ws.Range("I18").Value = 0.6198148578431102
$ws.Range("L18").Value = 0.2127765233942398
$ws.Range("M18").Value = 0.2387502348259147
$ws.Range("N18").Value = 1.270571292514326
$ws.Range("O18").Value = 2.400368735230103
$ws.Range("B19").Value = 1.063313829966091
$ws.Range("C19").Value = 0.2360185571340878
$ws.Range("E19").Value = 0.106287788169297
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.0024196693664835
$ws.Range("I19").Value = 0.6202974153492669
$ws.Range("L19").Value = 0.2124785685008277
$ws.Range("M19").Value = 0.2375639271657946
$ws.Range("N19").Value = 1.271156419462812
$ws.Range("O19").Value = 2.400586133212158
$ws.Range("B20").Value = 1.095375953952896
$ws.Range("C20").Value = 0.2405012802727242
$ws.Range("E20").Value = 0.1059794923638435
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.00241879396917367
$ws.Range("I20").Value = 0.6181516362639776
$ws.Range("L20").Value = 0.2138260888517607
$ws.Range("M20").Value = 0.2429060365687548
$ws.Range("N20").Value = 1.268555934841942
$ws.Range("O20").Value = 2.399717910063799
$ws.Range("B21").Value = 1.203165832043737
$ws.Range("C21").Value = 0.255490819022782
$ws.Range("E21").Value = 0.1049851005024554
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.002415948567998071
$ws.Range("I21").Value = 0.6113937154807623
$ws.Range("L21").Value = 0.2184703723113586
$ws.Range("M21").Value = 0.2609166923254094
$ws.Range("N21").Value = 1.260389563962477
$ws.Range("O21").Value = 2.39872924842075
$ws.Range("B22").Value = 1.273654374632258
$ws.Range("C22").Value = 0.2652330449294311
$ws.Range("E22").Value = 0.1043660430787217
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002414160056544225
$ws.Range("I22").Value = 0.6073140830137795
$ws.Range("L22").Value = 0.2215923941707985
$ws.Range("M22").Value = 0.2727326497463451
$ws.Range("N22").Value = 1.255477266517886
$ws.Range("O22").Value = 2.399522610089349
$ws.Range("B23").Value = 1.236029967128275
$ws.Range("C23").Value = 0.2600384553234676
$ws.Range("E23").Value = 0.1046936274952666
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.002415108129327538
$ws.Range("I23").Value = 0.6094605883122455
$ws.Range("L23").Value = 0.2199181997372932
$ws.Range("M23").Value = 0.2664222135953764
$ws.Range("N23").Value = 1.258060211958259
$ws.Range("O23").Value = 2.398967321383424
$ws.Range("B24").Value = 1.093614990347874
$ws.Range("C24").Value = 0.2402553786168653
$ws.Range("E24").Value = 0.1059962664696661
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002418841680836887
$ws.Range("I24").Value = 0.6182677744598877
$ws.Range("L24").Value = 0.2137516492075804
$ws.Range("M24").Value = 0.2426124375608438
$ws.Range("N24").Value = 1.268696591718992
$ws.Range("O24").Value = 2.399758339479945
$ws.Range("B25").Value = 0.9402223319801806
$ws.Range("C25").Value = 0.2186880470454469
$ws.Range("E25").Value = 0.1075340914646712
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.002423176157997747
$ws.Range("I25").Value = 0.6292118113341054
$ws.Range("L25").Value = 0.2074744823463348
$ws.Range("M25").Value = 0.2171303047389728
$ws.Range("N25").Value = 1.281997865483554
$ws.Range("O25").Value = 2.406777472351678
